# 3-mascaras-ejemplo.xlsx
# Adds a new subnet-mask example row pair (rows 7-8) to "Hoja1":
#   Row 7 = the merged "mask" display row (255.255.224.0  ->  /19)
#   Row 8 = the individual-bit row backing that mask
# and updates the sheet view (zoom + selection) to match where the
# author ended up after typing the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row 7: merged "decimal mask" cells (one value per merged block) ---
$ws.Range("A7").Value = 255
$ws.Range("I7").Value = 255
$ws.Range("Q7").Value = 224
$ws.Range("Y7").Value = 0

# --- Row 8: the 32 individual mask bits -------------------------------
# Octet 1 (A:H) = 255 -> 1 1 1 1 1 1 1 1
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1

# Octet 2 (I:P) = 255 -> 1 1 1 1 1 1 1 1
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 1
$ws.Range("P8").Value = 1

# Octet 3 (Q:X) = 224 -> 1 1 1 0 0 0 0 0
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 1
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 0
$ws.Range("U8").Value = 0
$ws.Range("V8").Value = 0
$ws.Range("W8").Value = 0
$ws.Range("X8").Value = 0

# Octet 4 (Y:AF) = 0 -> 0 0 0 0 0 0 0 0
$ws.Range("Y8").Value = 0
$ws.Range("Z8").Value = 0
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 0
$ws.Range("AC8").Value = 0
$ws.Range("AD8").Value = 0
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 0

# --- View state: the author zoomed out a bit and left the selection ---
# --- on the cell they had just typed into (T8). -----------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 160
$ws.Range("T8").Select()
